# Add a "time_taken" column (F) with per-row timestamps to the panel sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - same style as the other headers (B1:E1), which use style index 1
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Per-row timestamp values (as text, matching the inlineStr cells from the source)
$timestamps = @(
    "2021-10-05 10:52:06.014503",
    "2021-10-05 10:52:06.014514",
    "2021-10-05 10:52:06.014518",
    "2021-10-05 10:52:06.014521",
    "2021-10-05 10:52:06.014524",
    "2021-10-05 10:52:06.014527",
    "2021-10-05 10:52:06.014529",
    "2021-10-05 10:52:06.014532",
    "2021-10-05 10:52:06.014534",
    "2021-10-05 10:52:06.014537",
    "2021-10-05 10:52:06.014539",
    "2021-10-05 10:52:06.014542",
    "2021-10-05 10:52:06.014545",
    "2021-10-05 10:52:06.014547",
    "2021-10-05 10:52:06.014550",
    "2021-10-05 10:52:06.014553",
    "2021-10-05 10:52:06.014556",
    "2021-10-05 10:52:06.014559",
    "2021-10-05 10:52:06.014561",
    "2021-10-05 10:52:06.014564",
    "2021-10-05 10:52:06.014566",
    "2021-10-05 10:52:06.014569",
    "2021-10-05 10:52:06.014571",
    "2021-10-05 10:52:06.014574",
    "2021-10-05 10:52:06.014577",
    "2021-10-05 10:52:06.014580",
    "2021-10-05 10:52:06.014582",
    "2021-10-05 10:52:06.014585",
    "2021-10-05 10:52:06.014588",
    "2021-10-05 10:52:06.014590",
    "2021-10-05 10:52:06.014593",
    "2021-10-05 10:52:06.014596",
    "2021-10-05 10:52:06.014598",
    "2021-10-05 10:52:06.014601",
    "2021-10-05 10:52:06.014604",
    "2021-10-05 10:52:06.014607",
    "2021-10-05 10:52:06.014609",
    "2021-10-05 10:52:06.014612",
    "2021-10-05 10:52:06.014614",
    "2021-10-05 10:52:06.014617"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
